$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DT; this shifts the existing DT -> DU and DU -> DV
# (i.e. "nom" and "url_produit" move one column to the right), and also
# bumps the sheet dimension from DU206 to DV206 automatically.
$ws.Columns("DT:DT").Insert()

# New header cell for the freshly inserted column (row 1): a snapshot timestamp.
$ws.Range("DT1").Value = "2026-02-02 10:34:13"

# For the data rows that already have a numeric price in column DS (rows 2-80),
# mirror that same price into the newly inserted column DT.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 123).Value2   # column 123 = DS
    $ws.Cells.Item($r, 124).Value = $price     # column 124 = DT
}

# Rows 81-206 have no price recorded in DS (blank), so the corresponding new
# DT cell is left blank as well - nothing further to do for those rows.
